# Update the "want to go" (想去人数) counts in column F across the four
# sheets of the workbook, reflecting a refreshed scrape of bilibili
# exhibition/show data (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 960
$ws1.Range("F6").Value  = 5238
$ws1.Range("F17").Value = 1749
$ws1.Range("F18").Value = 1451
$ws1.Range("F19").Value = 820
$ws1.Range("F23").Value = 505
$ws1.Range("F25").Value = 1045
$ws1.Range("F28").Value = 2526
$ws1.Range("F31").Value = 50
$ws1.Range("F40").Value = 634

# Sheet "演出" (Shows/Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 8

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 230

# Sheet "全部类型" (All types - combined view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 230
$ws4.Range("F5").Value  = 960
$ws4.Range("F7").Value  = 5238
$ws4.Range("F21").Value = 8
$ws4.Range("F23").Value = 1749
$ws4.Range("F24").Value = 1451
$ws4.Range("F25").Value = 820
$ws4.Range("F29").Value = 505
$ws4.Range("F31").Value = 1045
$ws4.Range("F33").Value = 2526
